$d = $word.ActiveDocument

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that sits after " non vide"
#    (it is being relocated to the very end of the document below).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*non vide*") {
        $nonVideXml = '<w:p ' + $wns + '><w:pPr><w:rPr><w:i/><w:color w:val="2E74B5" w:themeColor="accent1" w:themeShade="BF"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:color w:val="2E74B5" w:themeColor="accent1" w:themeShade="BF"/></w:rPr><w:t>Malloc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/><w:color w:val="2E74B5" w:themeColor="accent1" w:themeShade="BF"/></w:rPr><w:t xml:space="preserve"> allocation de mémoire peut </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:color w:val="2E74B5" w:themeColor="accent1" w:themeShade="BF"/></w:rPr><w:t>etre</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/><w:color w:val="2E74B5" w:themeColor="accent1" w:themeShade="BF"/></w:rPr><w:t xml:space="preserve"> non vide</w:t></w:r></w:p>'
        $p.Range.InsertXML($nonVideXml)
        break
    }
}

# ------------------------------------------------------------------
# 2) Append the new "exo structure" paragraphs at the end of the
#    document, just before the lone empty paragraph that precedes
#    <w:sectPr> (that empty paragraph stays put, between "Heritage"
#    and "Argc").
# ------------------------------------------------------------------
$last = $d.Paragraphs.Last
$tailRange = $last.Range
$tailRange.InsertParagraphBefore()

$classeXml = '<w:p ' + $wns + '>' +
    '<w:r><w:t xml:space="preserve">Classe </w:t></w:r>' +
    '<w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>private</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'

$newParaIndex = $d.Paragraphs.Count - 1
$d.Paragraphs.Item($newParaIndex).Range.InsertXML($classeXml)

$heritageXml = '<w:p ' + $wns + '>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Heritage</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>protected</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'

$last2 = $d.Paragraphs.Last
$last2.Range.InsertParagraphBefore()
$newParaIndex2 = $d.Paragraphs.Count - 1
$d.Paragraphs.Item($newParaIndex2).Range.InsertXML($heritageXml)

# The existing lone empty paragraph (originally the very last one in
# the document) is now just before $d.Paragraphs.Last - leave it
# untouched, and insert the remaining two paragraphs after it.

$last3 = $d.Paragraphs.Last
$last3.Range.InsertParagraphBefore()

$argcXml = '<w:p ' + $wns + '>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Argc</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> argument count </w:t></w:r>' +
    '<w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>int</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'

$newParaIndex3 = $d.Paragraphs.Count - 1
$d.Paragraphs.Item($newParaIndex3).Range.InsertXML($argcXml)

$last4 = $d.Paragraphs.Last
$last4.Range.InsertParagraphBefore()

$argvXml = '<w:p ' + $wns + '>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Argv</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t> : double pointeur</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '</w:p>'

$newParaIndex4 = $d.Paragraphs.Count - 1
$d.Paragraphs.Item($newParaIndex4).Range.InsertXML($argvXml)

Write-Output "done"
